$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H11").Value = 265.21054
$ws_ALC.Range("I11").Value = 265.21054
$ws_ALC.Range("K11").Value = 265.21054
$ws_ALC.Range("M11").Value = -125.21054

$ws_ALC.Range("H69").Value = 5900
$ws_ALC.Range("J69").Value = 0
$ws_ALC.Range("L69").Value = 0
$ws_ALC.Range("N69").ClearContents()

$ws_ALC.Range("H72").Value = 5900
$ws_ALC.Range("J72").Value = 0
$ws_ALC.Range("L72").Value = 0
$ws_ALC.Range("N72").ClearContents()

$ws_ALC.Range("H132").Value = 7829.8125
$ws_ALC.Range("I132").Value = 8655.071
$ws_ALC.Range("K132").Value = 25965.213
$ws_ALC.Range("M132").Value = -23435.213

$ws_ALC.Range("H135").Value = 1357.8125
$ws_ALC.Range("I135").Value = 1441.6666
$ws_ALC.Range("J135").Value = 100
$ws_ALC.Range("K135").Value = 12974.9994
$ws_ALC.Range("L135").Value = 900
$ws_ALC.Range("M135").Value = -10439.9994
$ws_ALC.Range("N135").Value = -5970

$ws_ALC.Range("H137").Value = 1343.4688
$ws_ALC.Range("J137").Value = 1714.7
$ws_ALC.Range("L137").Value = 5144.1
$ws_ALC.Range("N137").Value = -10244.1

$ws_ALC.Range("H138").Value = 2768.3455
$ws_ALC.Range("I138").Value = 1878.76
$ws_ALC.Range("J138").Value = 3509.6667
$ws_ALC.Range("K138").Value = 5636.28
$ws_ALC.Range("L138").Value = 10529.0001
$ws_ALC.Range("M138").Value = -496.2799999999997
$ws_ALC.Range("N138").Value = -20809.0001

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 4099.684
$ws_ARM.Range("I32").Value = 4099.684
$ws_ARM.Range("J32").Value = 0
$ws_ARM.Range("K32").Value = 4099.684
$ws_ARM.Range("L32").Value = 0
$ws_ARM.Range("M32").Value = -3812.684
$ws_ARM.Range("N32").ClearContents()

$ws_ARM.Range("H45").Value = 97287.77
$ws_ARM.Range("I45").Value = 171194.75
$ws_ARM.Range("J45").Value = 8599.4
$ws_ARM.Range("K45").Value = 171194.75
$ws_ARM.Range("L45").Value = 8599.4
$ws_ARM.Range("M45").Value = -170817.75
$ws_ARM.Range("N45").Value = -9353.4

$ws_ARM.Range("H61").Value = 10423047
$ws_ARM.Range("I61").Value = 12200207
$ws_ARM.Range("K61").Value = 12200207
$ws_ARM.Range("M61").Value = -12199995

$ws_ARM.Range("H110").Value = 4686.7896
$ws_ARM.Range("I110").Value = 3466.3076
$ws_ARM.Range("J110").Value = 7331.1665
$ws_ARM.Range("K110").Value = 3466.3076
$ws_ARM.Range("L110").Value = 7331.1665
$ws_ARM.Range("M110").Value = -1421.3076
$ws_ARM.Range("N110").Value = -11421.1665

$ws_ARM.Range("H122").Value = 2500.2727
$ws_ARM.Range("I122").Value = 2216.6667
$ws_ARM.Range("K122").Value = 6650.000100000001
$ws_ARM.Range("M122").Value = -4200.000100000001

$ws_ARM.Range("H132").Value = 3108.5833
$ws_ARM.Range("I132").Value = 2569.4285
$ws_ARM.Range("J132").Value = 6882.6665
$ws_ARM.Range("K132").Value = 7708.2855
$ws_ARM.Range("L132").Value = 20647.9995
$ws_ARM.Range("M132").Value = -5178.2855
$ws_ARM.Range("N132").Value = -25707.9995

$ws_ARM.Range("H136").Value = 10423047
$ws_ARM.Range("I136").Value = 12200207
$ws_ARM.Range("K136").Value = 36600621
$ws_ARM.Range("M136").Value = -36598071

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H64").Value = 397.5
$ws_BSM.Range("J64").Value = 437
$ws_BSM.Range("L64").Value = 437
$ws_BSM.Range("N64").Value = -887

$ws_BSM.Range("H67").Value = 397.5
$ws_BSM.Range("J67").Value = 437
$ws_BSM.Range("L67").Value = 437
$ws_BSM.Range("N67").Value = -1997

$ws_BSM.Range("H99").Value = 5092.5
$ws_BSM.Range("I99").Value = 3998.889
$ws_BSM.Range("J99").Value = 6498.5713
$ws_BSM.Range("K99").Value = 3998.889
$ws_BSM.Range("L99").Value = 6498.5713
$ws_BSM.Range("M99").Value = -2500.889
$ws_BSM.Range("N99").Value = -9494.5713

$ws_BSM.Range("H105").Value = 1116.2222
$ws_BSM.Range("I105").Value = 1033.5
$ws_BSM.Range("K105").Value = 1033.5
$ws_BSM.Range("M105").Value = 713.5

$ws_BSM.Range("H134").Value = 2353.05
$ws_BSM.Range("I134").Value = 2377.487
$ws_BSM.Range("K134").Value = 7132.461
$ws_BSM.Range("M134").Value = -4597.461

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H4").Value = 2998.8333
$ws_CRP.Range("J4").Value = 3664.6667
$ws_CRP.Range("L4").Value = 3664.6667
$ws_CRP.Range("N4").Value = -3888.6667

$ws_CRP.Range("H7").Value = 95.57143000000001
$ws_CRP.Range("I7").Value = 95.57143000000001
$ws_CRP.Range("K7").Value = 95.57143000000001
$ws_CRP.Range("M7").Value = 17.42856999999999

$ws_CRP.Range("H31").Value = 5325.2256
$ws_CRP.Range("I31").Value = 4245
$ws_CRP.Range("J31").Value = 7035.5835
$ws_CRP.Range("K31").Value = 4245
$ws_CRP.Range("L31").Value = 7035.5835
$ws_CRP.Range("M31").Value = -3950
$ws_CRP.Range("N31").Value = -7625.5835

$ws_CRP.Range("H34").Value = 5325.2256
$ws_CRP.Range("I34").Value = 4245
$ws_CRP.Range("J34").Value = 7035.5835
$ws_CRP.Range("K34").Value = 4245
$ws_CRP.Range("L34").Value = 7035.5835
$ws_CRP.Range("M34").Value = -4043
$ws_CRP.Range("N34").Value = -7439.5835

$ws_CRP.Range("H99").Value = 7759.6
$ws_CRP.Range("I99").Value = 7966.3335
$ws_CRP.Range("J99").Value = 7449.5
$ws_CRP.Range("K99").Value = 7966.3335
$ws_CRP.Range("L99").Value = 7449.5
$ws_CRP.Range("M99").Value = -6468.3335
$ws_CRP.Range("N99").Value = -10445.5

$ws_CRP.Range("H107").Value = 460.8889
$ws_CRP.Range("I107").Value = 378.2857
$ws_CRP.Range("K107").Value = 378.2857
$ws_CRP.Range("M107").Value = 1541.7143

$ws_CRP.Range("H126").Value = 7759.6
$ws_CRP.Range("I126").Value = 7966.3335
$ws_CRP.Range("J126").Value = 7449.5
$ws_CRP.Range("K126").Value = 23899.0005
$ws_CRP.Range("L126").Value = 22348.5
$ws_CRP.Range("M126").Value = -21429.0005
$ws_CRP.Range("N126").Value = -27288.5

$ws_CRP.Range("H134").Value = 6305.684
$ws_CRP.Range("I134").Value = 3625.9167
$ws_CRP.Range("K134").Value = 10877.7501
$ws_CRP.Range("M134").Value = -8342.750100000001

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H8").Value = 854.55554
$ws_CUL.Range("I8").Value = 854.55554
$ws_CUL.Range("K8").Value = 2563.66662
$ws_CUL.Range("M8").Value = -2424.66662

$ws_CUL.Range("H12").Value = 266.08334
$ws_CUL.Range("I12").Value = 12.428572
$ws_CUL.Range("J12").Value = 621.2
$ws_CUL.Range("K12").Value = 37.285716
$ws_CUL.Range("L12").Value = 1863.6
$ws_CUL.Range("M12").Value = 135.714284
$ws_CUL.Range("N12").Value = -2209.6

$ws_CUL.Range("H62").Value = 0
$ws_CUL.Range("J62").Value = 0
$ws_CUL.Range("L62").Value = 0
$ws_CUL.Range("N62").ClearContents()

$ws_CUL.Range("H65").Value = 0
$ws_CUL.Range("J65").Value = 0
$ws_CUL.Range("L65").Value = 0
$ws_CUL.Range("N65").ClearContents()

$ws_CUL.Range("H97").Value = 1462.36
$ws_CUL.Range("I97").Value = 1271.8667
$ws_CUL.Range("J97").Value = 1748.1
$ws_CUL.Range("K97").Value = 3815.6001
$ws_CUL.Range("L97").Value = 5244.299999999999
$ws_CUL.Range("M97").Value = -3319.6001
$ws_CUL.Range("N97").Value = -6236.299999999999

$ws_CUL.Range("H107").Value = 933.1579
$ws_CUL.Range("J107").Value = 903.0833
$ws_CUL.Range("L107").Value = 2709.2499
$ws_CUL.Range("N107").Value = -6549.2499

$ws_CUL.Range("H109").Value = 2500
$ws_CUL.Range("I109").Value = 0
$ws_CUL.Range("K109").Value = 0
$ws_CUL.Range("M109").ClearContents()

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H113").Value = 20850.166
$ws_GSM.Range("J113").Value = 2333.3333
$ws_GSM.Range("L113").Value = 2333.3333
$ws_GSM.Range("N113").Value = -6673.3333

$ws_GSM.Range("H122").Value = 1973.4445
$ws_GSM.Range("I122").Value = 2037.5714
$ws_GSM.Range("K122").Value = 6112.7142
$ws_GSM.Range("M122").Value = -3662.7142

$ws_GSM.Range("H130").Value = 0
$ws_GSM.Range("I130").Value = 0
$ws_GSM.Range("K130").Value = 0
$ws_GSM.Range("M130").ClearContents()

$ws_GSM.Range("H132").Value = 2401.7778
$ws_GSM.Range("I132").Value = 2049.75
$ws_GSM.Range("J132").Value = 3105.8333
$ws_GSM.Range("K132").Value = 6149.25
$ws_GSM.Range("L132").Value = 9317.499899999999
$ws_GSM.Range("M132").Value = -3619.25
$ws_GSM.Range("N132").Value = -14377.4999

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H22").Value = 2141.5
$ws_LTW.Range("J22").Value = 2612.5
$ws_LTW.Range("L22").Value = 2612.5
$ws_LTW.Range("N22").Value = -3202.5

$ws_LTW.Range("H27").Value = 2141.5
$ws_LTW.Range("J27").Value = 2612.5
$ws_LTW.Range("L27").Value = 2612.5
$ws_LTW.Range("N27").Value = -2826.5

$ws_LTW.Range("H40").Value = 3330.85
$ws_LTW.Range("I40").Value = 3139.3684
$ws_LTW.Range("K40").Value = 3139.3684
$ws_LTW.Range("M40").Value = -3003.3684

$ws_LTW.Range("H61").Value = 170542.33
$ws_LTW.Range("I61").Value = 253126
$ws_LTW.Range("J61").Value = 5375
$ws_LTW.Range("K61").Value = 253126
$ws_LTW.Range("L61").Value = 5375
$ws_LTW.Range("M61").Value = -252924
$ws_LTW.Range("N61").Value = -5779

$ws_LTW.Range("H113").Value = 170542.33
$ws_LTW.Range("I113").Value = 253126
$ws_LTW.Range("J113").Value = 5375
$ws_LTW.Range("K113").Value = 253126
$ws_LTW.Range("L113").Value = 5375
$ws_LTW.Range("M113").Value = -250956
$ws_LTW.Range("N113").Value = -9715

$ws_LTW.Range("H132").Value = 9529.541999999999
$ws_LTW.Range("I132").Value = 11631.1875
$ws_LTW.Range("J132").Value = 5326.25
$ws_LTW.Range("K132").Value = 34893.5625
$ws_LTW.Range("L132").Value = 15978.75
$ws_LTW.Range("M132").Value = -32363.5625
$ws_LTW.Range("N132").Value = -21038.75

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H113").Value = 412.5
$ws_WVR.Range("I113").Value = 420.2381
$ws_WVR.Range("K113").Value = 1260.7143
$ws_WVR.Range("M113").Value = 909.2857000000001

$ws_WVR.Range("H122").Value = 3492.7144
$ws_WVR.Range("I122").Value = 1999.2354
$ws_WVR.Range("K122").Value = 5997.706200000001
$ws_WVR.Range("M122").Value = -3547.706200000001

$ws_WVR.Range("H126").Value = 1600.6666
$ws_WVR.Range("I126").Value = 1630.4286
$ws_WVR.Range("J126").Value = 1496.5
$ws_WVR.Range("K126").Value = 4891.2858
$ws_WVR.Range("L126").Value = 4489.5
$ws_WVR.Range("M126").Value = -2421.2858
$ws_WVR.Range("N126").Value = -9429.5

$ws_WVR.Range("H132").Value = 3968.224
$ws_WVR.Range("I132").Value = 3436.318
$ws_WVR.Range("K132").Value = 10308.954
$ws_WVR.Range("M132").Value = -7778.954000000002
